$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '26.716.07'
$ws.Cells.Item(2, 5).Value = '  -0.22%  '
$ws.Cells.Item(3, 4).Value = '1.633.81'
$ws.Cells.Item(3, 5).Value = '  -0.91%  '
$ws.Cells.Item(4, 5).Value = '  +0.06%  '
$ws.Cells.Item(5, 4).Value = '218.02'
$ws.Cells.Item(5, 5).Value = '  +0.63%  '
$ws.Cells.Item(6, 4).Value = '0.498'
$ws.Cells.Item(6, 5).Value = '  -1.57%  '
$ws.Cells.Item(7, 5).Value = '  +0.09%  '
$ws.Cells.Item(8, 5).Value = '  -1.42%  '
$ws.Cells.Item(9, 4).Value = '0.0620'
$ws.Cells.Item(9, 5).Value = '  -1.09%  '
$ws.Cells.Item(10, 4).Value = '18.98'
$ws.Cells.Item(10, 5).Value = '  -1.31%  '
$ws.Cells.Item(11, 4).Value = '0.0843'
$ws.Cells.Item(11, 5).Value = '  -0.17%  '
$ws.Cells.Item(12, 4).Value = '1.860.59'
$ws.Cells.Item(12, 5).Value = '  -0.95%  '
$ws.Cells.Item(13, 4).Value = '1.651.50'
$ws.Cells.Item(13, 5).Value = '  -0.86%  '
$ws.Cells.Item(14, 5).Value = '  -2.32%  '
$ws.Cells.Item(15, 5).Value = '  -2.19%  '
$ws.Cells.Item(16, 5).Value = '  -2.21%  '
$ws.Cells.Item(17, 4).Value = '26.679.55'
$ws.Cells.Item(17, 5).Value = '  -0.45%  '
$ws.Cells.Item(18, 4).Value = '0.0₃0720'
$ws.Cells.Item(18, 5).Value = '  -3.34%  '
$ws.Cells.Item(19, 5).Value = '  +0.05%  '
$ws.Cells.Item(20, 4).Value = '210.81'
$ws.Cells.Item(20, 5).Value = '  -3.52%  '
$ws.Cells.Item(21, 4).Value = '4.30'
$ws.Cells.Item(21, 5).Value = '  -1.87%  '
$ws.Cells.Item(22, 2).Value = 'Toncoin'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(22, 4).Value = '2.33'
$ws.Cells.Item(22, 5).Value = '  -8.25%  '
$ws.Cells.Item(23, 2).Value = 'Chainlink'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(23, 4).Value = '6.15'
$ws.Cells.Item(23, 5).Value = '  -2.51%  '
$ws.Cells.Item(24, 4).Value = '9.13'
$ws.Cells.Item(24, 5).Value = '  -3.56%  '
$ws.Cells.Item(25, 4).Value = '146.76'
$ws.Cells.Item(25, 5).Value = '  +0.52%  '
$ws.Cells.Item(26, 5).Value = '  +0.24%  '
$ws.Cells.Item(27, 5).Value = '  -2.57%  '
$ws.Cells.Item(28, 4).Value = '7.01'
$ws.Cells.Item(28, 5).Value = '  -2.61%  '
$ws.Cells.Item(29, 4).Value = '15.50'
$ws.Cells.Item(29, 5).Value = '  -2.03%  '
$ws.Cells.Item(30, 5).Value = '  -3.69%  '
$ws.Cells.Item(31, 5).Value = '  +0.55%  '
$ws.Cells.Item(32, 5).Value = '  -0.54%  '
$ws.Cells.Item(33, 5).Value = '  -3.05%  '
$ws.Cells.Item(34, 4).Value = '1.259.18'
$ws.Cells.Item(34, 5).Value = '  -1.74%  '
$ws.Cells.Item(35, 4).Value = '2.45'
$ws.Cells.Item(35, 5).Value = '  +0.29%  '
$ws.Cells.Item(36, 5).Value = '  -2.94%  '
$ws.Cells.Item(37, 5).Value = '  -3.93%  '
$ws.Cells.Item(38, 4).Value = '0.520'
$ws.Cells.Item(38, 5).Value = '  -4.36%  '
$ws.Cells.Item(39, 5).Value = '  +0.13%  '
$ws.Cells.Item(40, 4).Value = '0.797'
$ws.Cells.Item(40, 5).Value = '  -4.24%  '
$ws.Cells.Item(41, 4).Value = '0.798'
$ws.Cells.Item(41, 5).Value = '  -2.31%  '
$ws.Cells.Item(42, 5).Value = '  -4.12%  '
$ws.Cells.Item(43, 5).Value = '  -1.00%  '
$ws.Cells.Item(44, 4).Value = '5.25'
$ws.Cells.Item(44, 5).Value = '  -3.52%  '
$ws.Cells.Item(45, 4).Value = '90.91'
$ws.Cells.Item(45, 5).Value = '  -1.30%  '
$ws.Cells.Item(46, 4).Value = '59.71'
$ws.Cells.Item(46, 5).Value = '  -0.47%  '
$ws.Cells.Item(47, 5).Value = '  -2.73%  '
$ws.Cells.Item(48, 5).Value = '  +0.14%  '
$ws.Cells.Item(49, 2).Value = 'USDD'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Cells.Item(49, 4).Value = '1.01'
$ws.Cells.Item(49, 5).Value = '  +0.11%  '
$ws.Cells.Item(50, 2).Value = 'Mantle'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(50, 4).Value = '0.406'
$ws.Cells.Item(50, 5).Value = '  -0.51%  '
$ws.Cells.Item(51, 5).Value = '  -3.05%  '
